{"js": "// Insert the missing contact-info paragraph right after the name\n// heading (\"Dheeraj Chand\"), centered, matching the short-resume\n// contact line used in the long-form template.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the name paragraph (\"Dheeraj Chand\") \u2014 it is the first\n// paragraph in the document body for this template.\nlet namePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Dheeraj Chand\") {\n    namePara = paragraphs.items[i];\n    break;\n  }\n}\nif (!namePara) {\n  namePara = paragraphs.items[0];\n}\n\n// Insert a new, empty paragraph right after the name paragraph so the\n// new run does not inherit the bold/28pt run formatting used for the\n// name itself.\nconst contactPara = namePara.insertParagraph(\"\", Word.InsertLocation.after);\ncontactPara.alignment = Word.Alignment.centered;\n// Clear out any inherited paragraph-mark run formatting before adding\n// the plain (unformatted) contact-info text.\ncontactPara.clear();\ncontactPara.insertText(\n  \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Fix contact information missing from the short resume: add a\n# centered contact-info line directly under the name heading\n# (\"Dheeraj Chand\"), matching the long-resume contact line.\n\n$d = $word.ActiveDocument\n\n# Scope the search to the name paragraph only, so we never touch\n# anything else in the document even if \"Dheeraj Chand\" appears\n# again later in the body text.\n$namePara = $d.Paragraphs(1)\n$rng = $namePara.Range\n\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Dheeraj Chand\"\n# \"^p\" inserts a real paragraph break, producing a brand-new\n# (unformatted) paragraph right after the name for the contact\n# details, centered along with the name above it.\n$find.Replacement.Text = \"Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\"\n\n$find.Execute($Null, $False, $False, $False, $False, $False, $True, 1, $False, $Null, 1) | Out-Null\n\n# Center the newly created contact-info paragraph.\n$contactPara = $d.Paragraphs(2)\n$contactPara.Alignment = 1\n"}
